$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 14000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = ""
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3283.2246
$ws.Range("I61").Value = 1863.5
$ws.Range("K61").Value = 1863.5
$ws.Range("M61").Value = -1651.5
$ws.Range("H113").Value = 40932.5
$ws.Range("J113").Value = 40932.5
$ws.Range("L113").Value = 40932.5
$ws.Range("N113").Value = -49610.5
$ws.Range("H132").Value = 2201968
$ws.Range("I132").Value = 3911.7222
$ws.Range("J132").Value = 4529322
$ws.Range("K132").Value = 11735.1666
$ws.Range("L132").Value = 13587966
$ws.Range("M132").Value = -9205.1666
$ws.Range("N132").Value = -13593026
$ws.Range("H136").Value = 3283.2246
$ws.Range("I136").Value = 1863.5
$ws.Range("K136").Value = 5590.5
$ws.Range("M136").Value = -3040.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""
$ws.Range("H107").Value = 588.13043
$ws.Range("I107").Value = 350.0909
$ws.Range("J107").Value = 806.3333
$ws.Range("K107").Value = 350.0909
$ws.Range("L107").Value = 806.3333
$ws.Range("M107").Value = 1569.9091
$ws.Range("N107").Value = -4646.3333
$ws.Range("H132").Value = 2161.0715
$ws.Range("I132").Value = 1895.5883
$ws.Range("K132").Value = 5686.7649
$ws.Range("M132").Value = -3156.7649
$ws.Range("H134").Value = 11371055
$ws.Range("I134").Value = 17866210
$ws.Range("J134").Value = 4531.875
$ws.Range("K134").Value = 53598630
$ws.Range("L134").Value = 13595.625
$ws.Range("M134").Value = -53596095
$ws.Range("N134").Value = -18665.625
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1695.909
$ws.Range("I58").Value = 662
$ws.Range("K58").Value = 1986
$ws.Range("M58").Value = -1858
$ws.Range("H92").Value = 700
$ws.Range("I92").Value = 300
$ws.Range("K92").Value = 900
$ws.Range("M92").Value = 348
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 398.83334
$ws.Range("I107").Value = 403.91666
$ws.Range("J107").Value = 388.66666
$ws.Range("K107").Value = 403.91666
$ws.Range("L107").Value = 388.66666
$ws.Range("M107").Value = 1516.08334
$ws.Range("N107").Value = -4228.66666
$ws.Range("H132").Value = 2474.5
$ws.Range("I132").Value = 1700.8889
$ws.Range("J132").Value = 3469.1428
$ws.Range("K132").Value = 5102.6667
$ws.Range("L132").Value = 10407.4284
$ws.Range("M132").Value = -2572.6667
$ws.Range("N132").Value = -15467.4284
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5414.08
$ws.Range("I7").Value = 5397.4736
$ws.Range("K7").Value = 5397.4736
$ws.Range("M7").Value = -5285.4736
$ws.Range("H61").Value = 3776.5264
$ws.Range("I61").Value = 3268.1428
$ws.Range("J61").Value = 5200
$ws.Range("K61").Value = 3268.1428
$ws.Range("L61").Value = 5200
$ws.Range("M61").Value = -3066.1428
$ws.Range("N61").Value = -5604
$ws.Range("H93").Value = 9992.083000000001
$ws.Range("I93").Value = 13613.25
$ws.Range("J93").Value = 2749.75
$ws.Range("K93").Value = 13613.25
$ws.Range("L93").Value = 2749.75
$ws.Range("M93").Value = -12365.25
$ws.Range("N93").Value = -5245.75
$ws.Range("H113").Value = 3776.5264
$ws.Range("I113").Value = 3268.1428
$ws.Range("J113").Value = 5200
$ws.Range("K113").Value = 3268.1428
$ws.Range("L113").Value = 5200
$ws.Range("M113").Value = -1098.1428
$ws.Range("N113").Value = -9540
$ws.Range("H126").Value = 5414.08
$ws.Range("I126").Value = 5397.4736
$ws.Range("K126").Value = 16192.4208
$ws.Range("M126").Value = -13722.4208
$ws.Range("H132").Value = 3259.476
$ws.Range("I132").Value = 2976.6155
$ws.Range("J132").Value = 3719.125
$ws.Range("K132").Value = 8929.8465
$ws.Range("L132").Value = 11157.375
$ws.Range("M132").Value = -6399.8465
$ws.Range("N132").Value = -16217.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 91500
$ws.Range("J68").Value = 91500
$ws.Range("L68").Value = 91500
$ws.Range("N68").Value = -93122
$ws.Range("H71").Value = 91500
$ws.Range("J71").Value = 91500
$ws.Range("L71").Value = 274500
$ws.Range("N71").Value = -282612
$ws.Range("H96").Value = 5943
$ws.Range("I96").Value = 5276.6
$ws.Range("J96").Value = 6776
$ws.Range("K96").Value = 5276.6
$ws.Range("L96").Value = 6776
$ws.Range("M96").Value = -3903.6
$ws.Range("N96").Value = -9522
$ws.Range("H97").Value = 38967.43
$ws.Range("J97").Value = 38967.43
$ws.Range("L97").Value = 38967.43
$ws.Range("N97").Value = -40949.43
$ws.Range("H107").Value = 779.6667
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 779.6667
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2339.0001
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -6179.0001
$ws.Range("H113").Value = 995.2
$ws.Range("I113").Value = 1130.3334
$ws.Range("J113").Value = 792.5
$ws.Range("K113").Value = 3391.0002
$ws.Range("L113").Value = 2377.5
$ws.Range("M113").Value = -1221.0002
$ws.Range("N113").Value = -6717.5
$ws.Range("H122").Value = 1519.1428
$ws.Range("I122").Value = 1518.6471
$ws.Range("J122").Value = 1521.25
$ws.Range("K122").Value = 4555.9413
$ws.Range("L122").Value = 4563.75
$ws.Range("M122").Value = -2105.9413
$ws.Range("N122").Value = -9463.75
$ws.Range("H126").Value = 1182.7142
$ws.Range("I126").Value = 1196.7693
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 3590.3079
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -1120.3079
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 4863359
$ws.Range("I132").Value = 2056.7896
$ws.Range("J132").Value = 13260154
$ws.Range("K132").Value = 6170.3688
$ws.Range("L132").Value = 39780462
$ws.Range("M132").Value = -3640.3688
$ws.Range("N132").Value = -39785522
